$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell while forcing it to be stored as literal text,
# so numeric-looking strings (e.g. "1.00", "0.000274") are not coerced into numbers
# and original cell formatting/style is left untouched.
function Set-TextValue {
    param($CellRef, $Text)
    $rng = $ws.Range($CellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = $origStyle
}

Set-TextValue "D2" '69.605.69'
$ws.Range("E2").Value = '  +1.20%  '
Set-TextValue "D3" '3.890.00'
Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  +0.12%  '
Set-TextValue "D5" '605.28'
$ws.Range("E5").Value = '  +0.83%  '
Set-TextValue "D6" '170.24'
$ws.Range("E6").Value = '  +4.99%  '
Set-TextValue "D7" '3.888.01'
$ws.Range("E7").Value = '  +0.69%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +1.08%  '
Set-TextValue "D10" '0.168'
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("E12").Value = '  +2.27%  '
$ws.Range("E13").Value = '  +5.51%  '
Set-TextValue "D14" '38.35'
$ws.Range("E14").Value = '  +4.06%  '
Set-TextValue "D15" '4.548.16'
$ws.Range("E15").Value = '  +0.78%  '
Set-TextValue "D16" '3.881.47'
$ws.Range("E16").Value = '  +2.37%  '
Set-TextValue "D17" '69.634.62'
$ws.Range("E17").Value = '  +0.99%  '
Set-TextValue "D18" '18.72'
$ws.Range("E18").Value = '  +9.22%  '
Set-TextValue "D19" '7.65'
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("E20").Value = '  -0.71%  '
Set-TextValue "D21" '11.17'
$ws.Range("E21").Value = '  -1.93%  '
Set-TextValue "D22" '490.84'
$ws.Range("E22").Value = '  +1.53%  '
Set-TextValue "D23" '0.749'
$ws.Range("E23").Value = '  +4.41%  '
$ws.Range("E24").Value = '  +3.69%  '
$ws.Range("E25").Value = '  +1.78%  '
$ws.Range("E26").Value = '  +3.79%  '
Set-TextValue "D27" '12.36'
$ws.Range("E27").Value = '  +2.36%  '
Set-TextValue "D28" '10.18'
$ws.Range("E28").Value = '  +2.32%  '
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("E31").Value = '  +2.76%  '
Set-TextValue "D32" '4.043.25'
$ws.Range("E32").Value = '  +0.68%  '
Set-TextValue "D33" '7.82'
$ws.Range("E33").Value = '  -0.53%  '
Set-TextValue "D34" '31.97'
$ws.Range("E34").Value = '  -0.74%  '
Set-TextValue "D35" '3.856.57'
$ws.Range("E35").Value = '  +1.20%  '
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("E37").Value = '  +4.66%  '
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("E39").Value = '  +1.30%  '
Set-TextValue "D40" '3.29'
$ws.Range("E40").Value = '  +11.11%  '
Set-TextValue "D41" '1.00'
$ws.Range("E41").Value = '  +0.02%  '
Set-TextValue "D42" '0.329'
$ws.Range("E42").Value = '  +3.50%  '
Set-TextValue "D43" '2.10'
$ws.Range("E43").Value = '  +6.34%  '
Set-TextValue "D44" '437.12'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("E46").Value = '  +4.12%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("E48").Value = '  +3.32%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue "D49" '0.000274'
$ws.Range("E49").Value = '  +20.61%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D50" '143.85'
$ws.Range("E50").Value = '  +0.17%  '
Set-TextValue "D51" '40.21'
$ws.Range("E51").Value = '  +4.04%  '

Write-Output "Applied 81 cell updates to cryptos sheet."
